$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 2.763564666666667
$ws.Range("H2").Value = 8.290694
$ws.Range("I2").Value = 0.009909756414635561
$ws.Range("J2").Value = 0.009909756414635559
$ws.Range("M2").Value = 11.146846
$ws.Range("N2").Value = 33.440538
$ws.Range("O2").Value = 0.2594806085672136
$ws.Range("P2").Value = 0.2594806085672136
$ws.Range("Q2").Value = 30.80502975037468
$ws.Range("R2").Value = 277.245267753372
$ws.Range("S2").Value = 0.002571389625222483
$ws.Range("T2").Value = 0.002571389625222483

# Row 3
$ws.Range("G3").Value = 2.763564666666667
$ws.Range("H3").Value = 8.290694
$ws.Range("I3").Value = 0.009909756414635561
$ws.Range("J3").Value = 0.009909756414635559
$ws.Range("O3").Value = 0.6444737471070977
$ws.Range("P3").Value = 0.6444737471070977
$ws.Range("Q3").Value = 76.510661288306
$ws.Range("R3").Value = 688.595951594754
$ws.Range("S3").Value = 0.006386577849458778
$ws.Range("T3").Value = 0.006386577849458777

# Row 4
$ws.Range("G4").Value = 2.763564666666667
$ws.Range("H4").Value = 8.290694
$ws.Range("I4").Value = 0.009909756414635561
$ws.Range("J4").Value = 0.009909756414635559
$ws.Range("O4").Value = 0.09604564432568881
$ws.Range("P4").Value = 0.09604564432568881
$ws.Range("Q4").Value = 11.40235082376245
$ws.Range("R4").Value = 102.621157413862
$ws.Range("S4").Value = 0.0009517889399543002
$ws.Range("T4").Value = 0.0009517889399543

# Row 5
$ws.Range("G5").Value = 266.1315866666666
$ws.Range("H5").Value = 798.3947599999999
$ws.Range("I5").Value = 0.9543106517164206
$ws.Range("J5").Value = 0.9543106517164204
$ws.Range("M5").Value = 11.146846
$ws.Range("N5").Value = 33.440538
$ws.Range("O5").Value = 0.2594806085672136
$ws.Range("P5").Value = 0.2594806085672136
$ws.Range("Q5").Value = 2966.527812308987
$ws.Range("R5").Value = 26698.75031078088
$ws.Range("S5").Value = 0.247625108669551
$ws.Range("T5").Value = 0.247625108669551

# Row 6
$ws.Range("G6").Value = 266.1315866666666
$ws.Range("H6").Value = 798.3947599999999
$ws.Range("I6").Value = 0.9543106517164206
$ws.Range("J6").Value = 0.9543106517164204
$ws.Range("O6").Value = 0.6444737471070977
$ws.Range("P6").Value = 0.6444737471070977
$ws.Range("Q6").Value = 7367.985244265238
$ws.Range("R6").Value = 66311.86719838715
$ws.Range("S6").Value = 0.615028161615898
$ws.Range("T6").Value = 0.6150281616158979

# Row 7
$ws.Range("G7").Value = 266.1315866666666
$ws.Range("H7").Value = 798.3947599999999
$ws.Range("I7").Value = 0.9543106517164206
$ws.Range("J7").Value = 0.9543106517164204
$ws.Range("O7").Value = 0.09604564432568881
$ws.Range("P7").Value = 0.09604564432568881
$ws.Range("R7").Value = 9882.42894314548
$ws.Range("S7").Value = 0.09165738143097162
$ws.Range("T7").Value = 0.0916573814309716

# Row 8
$ws.Range("I8").Value = 0.03577959186894402
$ws.Range("J8").Value = 0.03577959186894401
$ws.Range("M8").Value = 11.146846
$ws.Range("N8").Value = 33.440538
$ws.Range("O8").Value = 0.2594806085672136
$ws.Range("P8").Value = 0.2594806085672136
$ws.Range("Q8").Value = 111.2228541108513
$ws.Range("R8").Value = 1001.005686997662
$ws.Range("S8").Value = 0.009284110272440121
$ws.Range("T8").Value = 0.009284110272440119

# Row 9
$ws.Range("I9").Value = 0.03577959186894402
$ws.Range("J9").Value = 0.03577959186894401
$ws.Range("O9").Value = 0.6444737471070977
$ws.Range("P9").Value = 0.6444737471070977
$ws.Range("S9").Value = 0.023059007641741
$ws.Range("T9").Value = 0.02305900764174099

# Row 10
$ws.Range("I10").Value = 0.03577959186894402
$ws.Range("J10").Value = 0.03577959186894401
$ws.Range("O10").Value = 0.09604564432568881
$ws.Range("P10").Value = 0.09604564432568881
$ws.Range("S10").Value = 0.003436473954762905
$ws.Range("T10").Value = 0.003436473954762904
